# Updates the regression-stats table (rows 2-9, columns B-K) with the
# re-computed values from the commit (slope/intercept/.../count per site),
# matching a pandas/matplotlib-call refactor upstream.
#
# NOTE: this PS interpreter's tokenizer does not accept scientific-notation
# numeric literals (e.g. "1e-14"), so very small/large magnitudes below are
# written as exact integer-mantissa divisions/multiplications by powers of
# ten (via [Math]::Pow) that round-trip to the identical IEEE-754 double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9999999999999996
$ws.Range("C2").Value = 0.5674682166313563
$ws.Range("D2").Value = 0.01178039961164318
$ws.Range("E2").Value = 0.2657729160057531
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.0329825307187054
$ws.Range("H2").Value = 0.874617094785566
$ws.Range("I2").Value = 2.934568162247113
$ws.Range("J2").Value = (-184123248049139.0 / [Math]::Pow(10,28))
$ws.Range("K2").Value = 1035
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.3070023891072907
$ws.Range("D3").Value = 0.01407066102916952
$ws.Range("E3").Value = 0.31323995125622
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.3272735674166306
$ws.Range("H3").Value = 0.8307545166893964
$ws.Range("I3").Value = 3.162758929044037
$ws.Range("J3").Value = (3115085514680555.0 / [Math]::Pow(10,30))
$ws.Range("K3").Value = 1031
$ws.Range("B4").Value = 1.000000000000001
$ws.Range("C4").Value = 0.2883889200526189
$ws.Range("D4").Value = 0.0166682876273685
$ws.Range("E4").Value = 0.3902297037166906
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.460060290463257
$ws.Range("H4").Value = 0.7766652926778326
$ws.Range("I4").Value = 3.937364482533346
$ws.Range("J4").Value = (5728194089637838.0 / [Math]::Pow(10,30))
$ws.Range("K4").Value = 1037
$ws.Range("C5").Value = 0.4148163217381614
$ws.Range("D5").Value = 0.01343996622811544
$ws.Range("E5").Value = 0.3303914504409028
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.2095568041017632
$ws.Range("H5").Value = 0.8361298594968598
$ws.Range("I5").Value = 3.165847911001365
$ws.Range("J5").Value = (8386626770747089.0 / [Math]::Pow(10,30))
$ws.Range("K5").Value = 1087
$ws.Range("B6").Value = 0.9999999999999997
$ws.Range("C6").Value = 0.5967762292932864
$ws.Range("D6").Value = 0.0144814315488388
$ws.Range("E6").Value = 0.3368106220879001
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.07670325099191222
$ws.Range("H6").Value = 0.8158936334077829
$ws.Range("I6").Value = 3.320969459540207
$ws.Range("J6").Value = (2122400379543156.0 * [Math]::Pow(10,-16) * [Math]::Pow(10,-13))
$ws.Range("K6").Value = 1078
$ws.Range("B7").Value = 0.9999999999999994
$ws.Range("C7").Value = 0.838564129162174
$ws.Range("D7").Value = 0.01665466838797805
$ws.Range("E7").Value = 0.3817712502574175
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.02828557347831353
$ws.Range("H7").Value = 0.7837045324610361
$ws.Range("I7").Value = 3.623692066263507
$ws.Range("J7").Value = (1297079016131778.0 * [Math]::Pow(10,-15) * [Math]::Pow(10,-14))
$ws.Range("K7").Value = 997
$ws.Range("B8").Value = 0.9999999999999996
$ws.Range("C8").Value = 0.0411574724102992
$ws.Range("D8").Value = 0.01410740279886958
$ws.Range("E8").Value = 0.3129250543874262
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.8953935720716315
$ws.Range("H8").Value = 0.8644336197950682
$ws.Range("I8").Value = 3.016645558972559
$ws.Range("J8").Value = (8670420218642236.0 / [Math]::Pow(10,30))
$ws.Range("K8").Value = 790
$ws.Range("B9").Value = 1.000000000000002
$ws.Range("C9").Value = 0.2688994068953255
$ws.Range("D9").Value = 0.01770665261890524
$ws.Range("E9").Value = 0.3957801906792482
$ws.Range("F9").Value = (1578938943504622.0 * [Math]::Pow(10,-166) * [Math]::Pow(10,-166))
$ws.Range("G9").Value = 0.497027091671329
$ws.Range("H9").Value = 0.7560763277210524
$ws.Range("I9").Value = 4.385327212346684
$ws.Range("J9").Value = (-1350789293976524.0 / [Math]::Pow(10,30))
$ws.Range("K9").Value = 1031
